$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.636.83"
$ws.Range("E2").Value = "  -4.58%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.678.95"
$ws.Range("E3").Value = "  +1.84%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.26"
$ws.Range("E5").Value = "  -0.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.95"
$ws.Range("E6").Value = "  -5.69%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.599"
$ws.Range("E7").Value = "  -2.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.584"
$ws.Range("E9").Value = "  -2.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.67"
$ws.Range("E10").Value = "  -3.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0851"
$ws.Range("E11").Value = "  -1.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.16"
$ws.Range("E12").Value = "  -2.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.095.98"
$ws.Range("E13").Value = "  +2.16%  "

$ws.Range("E14").Value = "  +1.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.685.11"
$ws.Range("E15").Value = "  +2.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.938"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "15.20"
$ws.Range("E17").Value = "  -0.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "45.777.48"
$ws.Range("E18").Value = "  -4.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000102"
$ws.Range("E19").Value = "  -0.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.87"
$ws.Range("E20").Value = "  +0.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.91"
$ws.Range("E21").Value = "  -2.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.10"
$ws.Range("E22").Value = "  +2.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "284.11"
$ws.Range("E23").Value = "  +3.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.07"
$ws.Range("E24").Value = "  -0.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.25"
$ws.Range("E25").Value = "  +0.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "30.89"
$ws.Range("E26").Value = "  +0.60%  "

$ws.Range("E27").Value = "  -0.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.04"
$ws.Range("E28").Value = "  -1.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.68"
$ws.Range("E29").Value = "  -0.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.60"
$ws.Range("E30").Value = "  -3.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.22"
$ws.Range("E31").Value = "  -4.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.24"
$ws.Range("E32").Value = "  +0.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.77"
$ws.Range("E33").Value = "  +1.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.36"
$ws.Range("E34").Value = "  +6.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0848"
$ws.Range("E35").Value = "  -1.37%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "155.10"
$ws.Range("E36").Value = "  +2.08%  "

$ws.Range("E37").Value = "  -2.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.121"
$ws.Range("E38").Value = "  -2.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.85"
$ws.Range("E39").Value = "  +13.56%  "

$ws.Range("E40").Value = "  -0.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "16.34"
$ws.Range("E41").Value = "  -0.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.65"
$ws.Range("E42").Value = "  -1.04%  "

$ws.Range("E43").Value = "  -1.57%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.01"
$ws.Range("E44").Value = "  -6.88%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.113.38"
$ws.Range("E45").Value = "  -4.14%  "

$ws.Range("E46").Value = "  +0.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "94.28"
$ws.Range("E47").Value = "  -2.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "112.21"
$ws.Range("E48").Value = "  -2.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.33"
$ws.Range("E49").Value = "  -6.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.939.69"
$ws.Range("E50").Value = "  +2.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.201"
$ws.Range("E51").Value = "  -1.60%  "
